$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12882
$ws1.Range("F5").Value = 87
$ws1.Range("F10").Value = 12801
$ws1.Range("F13").Value = 8671
$ws1.Range("F14").Value = 7666
$ws1.Range("F15").Value = 191
$ws1.Range("F16").Value = 102
$ws1.Range("F19").Value = 980
$ws1.Range("F22").Value = 380
$ws1.Range("F23").Value = 183

# Sheet "全部类型" (All types) - same events, rows shifted by 1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12882
$ws4.Range("F6").Value = 87
$ws4.Range("F11").Value = 12801
$ws4.Range("F14").Value = 8671
$ws4.Range("F15").Value = 7666
$ws4.Range("F16").Value = 191
$ws4.Range("F17").Value = 102
$ws4.Range("F20").Value = 980
$ws4.Range("F24").Value = 380
$ws4.Range("F25").Value = 183
